# DB connection and Retry logic
# Adds a "UserId" header (renamed from "UserName") and introduces a third
# column holding the expected post-login welcome message ("Hello, RCG"),
# while the old "UserName" header slides over to the new column C.
# Applied identically to both the ParentCredentials and StudentCredentials
# sheets.

$wb = $excel.ActiveWorkbook

$originalActiveSheet = $wb.ActiveSheet

$sheetNames = @("ParentCredentials", "StudentCredentials")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    # Column A header becomes "UserId" (keeps its existing "Normal 2" cell
    # style / 12pt font - only the text changes).
    $ws.Range("A1").Value = "UserId"

    # The old "UserName" header text now lives in the new column C.
    $ws.Range("C1").Value = "UserName"

    # New column C, row 2: the expected greeting shown after a successful
    # login / retry, styled the same as the rest of the data row.
    $ws.Range("C2").Value = "Hello, RCG"
    $ws.Range("C2").Style = "Normal 2"

    # Column C should match column B's width (bestFit data column).
    $ws.Columns.Item(3).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

    # Move the sheet's selection to the newly added cell.
    $ws.Activate()
    $ws.Range("C2").Select()
}

# Restore whichever sheet was active before we started touching things.
$originalActiveSheet.Activate()
